$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(765).Insert()
$ws.Rows.Item(765).Insert()
$ws.Range("A765").Value = "2026/02/03"
$ws.Range("B765").Value = "火"
$ws.Range("C765").Value = 19
$ws.Range("D765").Value = 201
$ws.Range("A766").Value = "2026/02/03"
$ws.Range("B766").Value = "火"
$ws.Range("C766").Value = 22
$ws.Range("D766").Value = 201
